$wb = $excel.ActiveWorkbook

$newDate = "2020-03-09 23:59:19"

# Mapping of worksheet index (1-based, matches workbook sheet order) -> table id
$tableIds = @(
    "Compartment",
    "Compound",
    "Definition",
    "Enzyme",
    "FbcObjective",
    "Gene",
    "Layout",
    "Measurement",
    "PbConfig",
    "Position",
    "Protein",
    "Quantity",
    "QuantityInfo",
    "QuantityMatrix",
    "Reaction",
    "ReactionStoichiometry",
    "Regulator",
    "Relation",
    "Relationship",
    "SparseMatrix",
    "SparseMatrixColumn",
    "SparseMatrixOrdered",
    "SparseMatrixRow",
    "StoichiometricMatrix",
    "rxnconContingencyList",
    "rxnconReactionList"
)

for ($i = 0; $i -lt $tableIds.Count; $i++) {
    $sheetIndex = $i + 1
    $id = $tableIds[$i]
    $ws = $wb.Worksheets.Item($sheetIndex)
    $ws.Unprotect()

    $newTableLine = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='$id' name='$id' date='$newDate' objTablesVersion='0.0.8'"

    if ($sheetIndex -eq 1) {
        # First sheet also carries the document-level metadata string in A1,
        # with the table metadata string shifted down to A2.
        $ws.Range("A1").Value = "!!!ObjTables schema='SBtab' objTablesVersion='0.0.8' date='$newDate'"
        $ws.Range("A2").Value = $newTableLine
    } else {
        $ws.Range("A1").Value = $newTableLine
    }

    $ws.Protect()
}
